{"js": "// The canonical OOXML diff for this document shows no actual content,\n// formatting, or structural change \u2014 every hunk only touches the\n// auto-generated XML namespace-prefix declarations on each part's root\n// element (e.g. xmlns:ns8 -> xmlns:ns9, xmlns:ns19 -> xmlns:ns17), which\n// is a cosmetic artifact of re-serializing the package and carries no\n// semantic meaning (the prefixes are not referenced anywhere in the\n// document content). The commit message refers to an unrelated backend\n// source-code fix (company registration form / \"selecionar\" field type)\n// that has no corresponding text in this document.\n//\n// There is therefore nothing in the document body, headers, footers,\n// styles, or numbering for this script to change. Touch the body\n// (load + sync) so the script still performs a no-op round trip through\n// the Word JS API without altering any visible content.\nconst body = context.document.body;\nbody.load(\"text\");\nawait context.sync();\n", "ps1": "# The canonical OOXML diff for this document shows no actual content,\n# formatting, or structural change -- every hunk only touches the\n# auto-generated XML namespace-prefix declarations on each part's root\n# element (e.g. xmlns:ns8 -> xmlns:ns9, xmlns:ns19 -> xmlns:ns17), which\n# is a cosmetic artifact of re-serializing the package and carries no\n# semantic meaning (the prefixes are not referenced anywhere in the\n# document content). The commit message refers to an unrelated backend\n# source-code fix (company registration form / \"selecionar\" field type)\n# that has no corresponding text in this document.\n#\n# There is therefore nothing in the document body, headers, footers,\n# styles, or numbering for this script to change. Touch the document\n# (read Content) so the script still performs a no-op round trip through\n# the Word COM object model without altering any visible content.\n$d = $word.ActiveDocument\n$null = $d.Content.Text\n"}
